# Update Data Category 1
# - "Info" sheet: merge the "Derivation"/"Date of Publication of Original
#   Dataset" columns (F/G) into a single "Date of Publication" column (F),
#   and merge the "Temporal aggregation..." / publication-date columns (F/G)
#   in the data row into a single date value in column F. Columns G become
#   empty (fully cleared, not just blanked).
# - Active tab switches from "Selected" to "Info", with a new selection.

$wb = $excel.ActiveWorkbook
$wsSelected = $wb.Worksheets.Item("Selected")
$wsInfo = $wb.Worksheets.Item("Info")

# Header row: replace "Derivation" (F5) + "Date of Publication of Original
# Dataset" (G5) with a single "Date of Publication" header in F5; drop G5.
$wsInfo.Range("F5").Value = "Date of Publication"
$wsInfo.Range("G5").Clear()

# Data row: column F held descriptive text, column G held the actual date
# serial (45537). Keep the date in F, drop the old text and G entirely.
$wsInfo.Range("F6").Value = $wsInfo.Range("G6").Value()
$wsInfo.Range("G6").Clear()

# Make "Info" the active/selected tab with a fresh selection, and drop the
# previous tab-selected state from "Selected" (selection there is unchanged).
$wsInfo.Range("F16").Select() | Out-Null
$wsSelected.Activate() | Out-Null
$wsInfo.Activate() | Out-Null
